# sm_car_data_Tire_Tire2x.xlsx — "Update 2p0. Convention change to support
# multi-axle vehicles": the single "Bus_Makhulu_2x" tire-pair sheet is split
# into two sheets, one per axle/tire spec, each keeping the same layout.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the existing sheet so the new one inherits all formatting,
# conditional formatting, column widths, styles, etc.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the sheets to the new per-tire convention.
$ws1.Name = "Tire2x_270_70R22"
$ws2.Name = "Tire2x_430_50R38"

# --- Sheet 2: Tire2x_430_50R38 -------------------------------------------
$ws2.Range("H2").Value = "Tire"
$ws2.Range("H3").Value = "Tire2x_430_50R38"
# H7 keeps the original formula/value (0.1359705*2) inherited from the copy.

# --- Sheet 1: Tire2x_270_70R22 -------------------------------------------
# "Instance" (H3) now matches the sheet/tire name; "Type" (H2) stays "Tire".
$ws1.Range("H2").Value = "Tire"
$ws1.Range("H3").Value = "Tire2x_270_70R22"
# xOffset is now a plain authored value instead of a formula.
$ws1.Range("H7").Value = 0.4572

# Sheet 1 keeps all three conditional-formatting rules (already duplicated
# by Copy) and is left as the non-selected tab.

# Sheet 2 only keeps the class-highlight rule for the data block; the two
# extra rules that used to flag A19/A20 on the old single sheet are removed.
$ws2.Range("A20").FormatConditions.Delete()
$ws2.Range("A19").FormatConditions.Delete()

# Selection / active-tab bookkeeping to match the saved view state.
$ws1.Range("C25").Select()
$ws2.Activate()
$ws2.Range("J16").Select()
